$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (should be row 194 originally)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$species = @(
    ,("M_alsense", 103, 5189)
    ,("M_angelicum", 202, 5415)
    ,("M_diernhoferi", 141, 5603)
    ,("M_engbaekii", 380, 4098)
    ,("M_europaeum", 2674, 5189)
    ,("M_fragae", 36, 4441)
    ,("M_heraklionense", 418, 4568)
    ,("M_iranicum", 120, 5965)
    ,("M_paraense", 205, 5254)
    ,("M_persicum", 106, 5346)
    ,("M_szulgai", 240, 5697)
    ,("M_triviale", 10, 3416)
)

$r = $lastRow
foreach ($item in $species) {
    $r = $r + 1
    $name = $item[0]
    $hgt = $item[1]
    $total = $item[2]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $hgt
    $ws.Cells.Item($r, 3).Value = $total
    $ws.Cells.Item($r, 4).Formula = "=B" + $r + "/C" + $r + "*100"
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0.00"
}

$win = $excel.ActiveWindow
$win.ScrollRow = 177
$win.ScrollColumn = 1
$ws.Range("A198").Select()

Write-Output "Added rows up to $r"
